$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0000119090754144846
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 19575605.8673771
$ws.Range("E2").Value = 2459690191846.092
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2459709767453.586
